# New month cleanup / inventory additions
# Append 4 new inventory/contract rows (113-116) to the bottom of Sheet1,
# matching the existing "plain text" layout used by the rest of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$startRow = 113
$endRow   = 116

# Force the new cells to be stored as plain text (like every other data
# row in the sheet) instead of letting Excel auto-detect dates/numbers.
$ws.Range("A113:J116").NumberFormat = "@"

$newRows = @(
    @("4/30/2019", "SPE7L0-19-V-5061", "3",  "`$7,539.00 ", "6130014243614", "POWER SUPPLY",        "GEMS",    "144875",            "CP",  "2019 OCT 07"),
    @("4/30/2019", "SPE7MC-19-V-7482", "15", "`$4,344.00",  "5935014489052", "BACKSHELL,ELECTRICA", "Glenair", "712FS277NF1012-31", "CP",  "2019 OCT 07"),
    @("4/30/2019", "SPE5E3-19-V-7870", "3",  "`$391.14 ",   "5330013722277", "GASKET",              "Timken",  "01070-0279LOC6-11", "M33", "2019 OCT 07"),
    @("4/30/2019", "SPE5E3-19-V-7833", "3",  "`$391.41",    "5330013722276", "GASKET",              "Timken",  "01070-0279LOC6-9",  "M33", "2019 OCT 07")
)

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]
    for ($c = 1; $c -le $rowData.Count; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowData[$c - 1]
    }
}
